# Updated cryptos list on Sun Nov 17 05:14:37 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '90.810.13'
$ws.Range("E2").Value = '  -0.61%  '

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.124.13'
$ws.Range("E3").Value = '  -0.18%  '

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.21%  '

# Row 5 - Solana
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.04'
$ws.Range("E5").Value = '  +8.02%  '

# Row 6 - BNB
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '632.98'
$ws.Range("E6").Value = '  +1.46%  '

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.05'
$ws.Range("E7").Value = '  +8.29%  '

# Row 8 - Dogecoin
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.356'
$ws.Range("E8").Value = '  -5.63%  '

# Row 9 - USDC
$ws.Range("E9").Value = '  +0.03%  '

# Row 10 - LidoStakedEther
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.381.74'
$ws.Range("E10").Value = '  +8.15%  '

# Row 11 - Cardano
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.717'
$ws.Range("E11").Value = '  -4.75%  '

# Row 12 - TRON
$ws.Range("E12").Value = '  +3.75%  '

# Row 13 - Avalanche
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '36.30'
$ws.Range("E13").Value = '  +4.33%  '

# Row 14 - ShibaInu
$ws.Range("E14").Value = '  -3.02%  '

# Row 15 - Toncoin
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.60'
$ws.Range("E15").Value = '  +2.75%  '

# Row 16 - WrappedBTC
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '90.356.39'
$ws.Range("E16").Value = '  -0.81%  '

# Row 17 - WrappedliquidstakedEther2.0
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.694.33'
$ws.Range("E17").Value = '  -0.31%  '

# Row 18 - WrappedEther
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.111.18'

# Row 19 - SuiNetwork
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.70'
$ws.Range("E19").Value = '  -1.62%  '

# Row 20 - Chainlink->PEPE
$ws.Range("B20").Value = 'PEPE'
$ws.Range("C20").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000213'
$ws.Range("E20").Value = '  -5.47%  '

# Row 21 - PEPE->Chainlink
$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.32'
$ws.Range("E21").Value = '  +1.06%  '

# Row 22 - BitcoinCash
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '450.80'
$ws.Range("E22").Value = '  +3.89%  '

# Row 23 - Uniswap
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.00'
$ws.Range("E23").Value = '  +3.10%  '

# Row 24 - Polkadot
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.38'
$ws.Range("E24").Value = '  +4.05%  '

# Row 25 - NEARProtocol
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.87'
$ws.Range("E25").Value = '  -3.38%  '

# Row 26 - Aptos->Litecoin
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '90.49'
$ws.Range("E26").Value = '  +6.07%  '

# Row 27 - Litecoin->Aptos
$ws.Range("B27").Value = 'Aptos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.28'
$ws.Range("E27").Value = '  -0.54%  '

# Row 28 - WrappedeETH
$ws.Range("E28").Value = '  -0.64%  '

# Row 29 - Dai
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.997'
$ws.Range("E29").Value = '  -0.22%  '

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.28'
$ws.Range("E30").Value = '  +3.20%  '

# Row 31 - Cronos
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.160'
$ws.Range("E31").Value = '  -5.47%  '

# Row 32 - Binance-PegBSC-USD
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.993'
$ws.Range("E32").Value = '  +11.96%  '

# Row 33 - EthereumClassic
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.51'
$ws.Range("E33").Value = '  +16.34%  '

# Row 34 - Stellar
$ws.Range("E34").Value = '  +27.44%  '

# Row 35 - Kaspa
$ws.Range("E35").Value = '  +4.97%  '

# Row 36 - Bittensor
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '507.30'
$ws.Range("E36").Value = '  -4.67%  '

# Row 37 - PancakeSwap->dogwifhat
$ws.Range("B37").Value = 'dogwifhat'
$ws.Range("C37").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.75'
$ws.Range("E37").Value = '  -2.38%  '

# Row 38 - dogwifhat->PancakeSwap
$ws.Range("B38").Value = 'PancakeSwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.93'
$ws.Range("E38").Value = '  +4.22%  '

# Row 39 - RenderToken
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.04'
$ws.Range("E39").Value = '  -2.73%  '

# Row 40 - Fetch.AI
$ws.Range("E40").Value = '  +0.41%  '

# Row 41 - PolygonEcosystemToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.426'
$ws.Range("E41").Value = '  +12.26%  '

# Row 42 - WhiteBITCoin
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.20'
$ws.Range("E42").Value = '  -0.49%  '

# Row 43 - USDe
$ws.Range("E43").Value = '  +0.00%  '

# Row 44 - Hedera
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0847'
$ws.Range("E44").Value = '  +9.73%  '

# Row 45 - MantraDAO
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.23'
$ws.Range("E45").Value = '  +31.75%  '

# Row 46 - Stacks
$ws.Range("E46").Value = '  +2.07%  '

# Row 47 - ARBITRUM
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.694'
$ws.Range("E47").Value = '  +11.43%  '

# Row 48 - Monero
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '149.07'
$ws.Range("E48").Value = '  +2.86%  '

# Row 49 - Filecoin
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.53'
$ws.Range("E49").Value = '  +7.67%  '

# Row 50 - ImmutableX
$ws.Range("E50").Value = '  +3.99%  '

# Row 51 - OKB
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '44.83'
$ws.Range("E51").Value = '  +1.35%  '
